$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calibration data sorted ascending by time (column A), rows 2-8
$data = @(
    @(53229.207007, [double]"-3.6007941643e-05", [double]"-2.8861359878e-05", [double]"-7.6919675216e-06"),
    @(53240.071007, [double]"-0.00019395772465", [double]"-0.00015621063552", [double]"-4.4064187327e-05"),
    @(53252.939008, [double]"-0.00045995841893", [double]"-0.0003713917175", [double]"-9.7912477987e-05"),
    @(53266.939009, [double]"-0.0007229895", [double]"-0.0005913059", [double]"-0.0001485774"),
    @(53277.407009, [double]"-0.0004553877", [double]"-0.0003694314", [double]"-0.0001035892"),
    @(53288.67101, [double]"-0.00018249611083", [double]"-0.00014744250269", [double]"-4.3359405461e-05"),
    @(53299.07101, [double]"-2.9746389323e-05", [double]"-2.3840005802e-05", [double]"-6.3411136321e-06")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $row++
}
